$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Beta) values that changed
$ws.Range("C2").Value = 19.12075701903682
$ws.Range("E2").Value = 0.01982943797740053
$ws.Range("F2").Value = 9.544208811213233
$ws.Range("G2").Value = 9.155795358156086
$ws.Range("H2").Value = 9.924786459079462
$ws.Range("I2").Value = 0.002215635735467799
$ws.Range("J2").Value = 0.001967128511985355
$ws.Range("K2").Value = 0.002520713659240183
$ws.Range("L2").Value = 0.008320486266175504
$ws.Range("M2").Value = 0.008095330454988101
$ws.Range("N2").Value = 0.008549672873660905

# Update existing row 3 (Gamma) values that changed
$ws.Range("C3").Value = 0.04981522627320694
$ws.Range("D3").Value = 0.04815098319456564
$ws.Range("E3").Value = 0.0499839736740351
$ws.Range("F3").Value = 0.04688930007883621
$ws.Range("G3").Value = 0.04661235885223579
$ws.Range("H3").Value = 0.04715788661824954
$ws.Range("I3").Value = 0.04532134925438309
$ws.Range("J3").Value = 0.04505468904984604
$ws.Range("K3").Value = 0.04557890687093379
$ws.Range("L3").Value = 0.04694098153863968
$ws.Range("M3").Value = 0.04666407978586447
$ws.Range("N3").Value = 0.04720960807018319

# Add new row 4 (Beta + Gamma)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 19.17057224531003
$ws.Range("D4").Value = 0.05509537454402212
$ws.Range("E4").Value = 0.06981341165143562
$ws.Range("F4").Value = 9.591098111292069
$ws.Range("G4").Value = 9.202407717008322
$ws.Range("H4").Value = 9.971944345697711
$ws.Range("I4").Value = 0.04753698498985089
$ws.Range("J4").Value = 0.0470218175618314
$ws.Range("K4").Value = 0.04809962053017398
$ws.Range("L4").Value = 0.05526146780481519
$ws.Range("M4").Value = 0.05475941024085258
$ws.Range("N4").Value = 0.05575928094384409

# Match the formatting of column A data cells (bold, centered, thin border)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
